$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row Right count
$ws.Range("B11").Value = 5

# Update "Total" row Right count
$ws.Range("B12").Value = 45

# Update "Total" correct/total marks summary text
$ws.Range("E12").Value = "45/140"
